$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.227.00"
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("D3").Value = "3.486.32"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.15"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.83"
$ws.Range("E6").Value = "  +4.21%  "
$ws.Range("E7").Value = "  +13.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "3.487.52"
$ws.Range("E9").Value = "  +2.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.26"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("E11").Value = "  +2.66%  "
$ws.Range("E12").Value = "  +3.60%  "
$ws.Range("D13").Value = "4.087.86"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000193"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.73"
$ws.Range("E16").Value = "  +5.46%  "
$ws.Range("D17").Value = "65.242.89"
$ws.Range("E17").Value = "  +2.67%  "
$ws.Range("D18").Value = "3.479.35"
$ws.Range("E18").Value = "  +2.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.49"
$ws.Range("E19").Value = "  +3.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.43"
$ws.Range("E20").Value = "  +2.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "384.16"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.23"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.554"
$ws.Range("E23").Value = "  +4.47%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.70"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000120"
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.05"
$ws.Range("E27").Value = "  +6.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.180"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  +12.43%  "
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  +1.95%  "
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("E34").Value = "  +5.85%  "
$ws.Range("E35").Value = "  +13.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.22"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("E37").Value = "  +5.74%  "
$ws.Range("D38").Value = "3.020.14"
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0781"
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.85"
$ws.Range("E40").Value = "  +6.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.94"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0323"
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.57"
$ws.Range("E43").Value = "  +5.33%  "
$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.783"
$ws.Range("E45").Value = "  +2.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.07"
$ws.Range("E46").Value = "  +11.53%  "
$ws.Range("E47").Value = "  +3.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "319.91"
$ws.Range("E48").Value = "  +9.19%  "
$ws.Range("E49").Value = "  +6.89%  "
$ws.Range("E50").Value = "  +5.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.878"
$ws.Range("E51").Value = "  +5.11%  "
